$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 23 - Complaint - Lockout No Access Users
$ws.Range("B23").Value = "Complaint – Lockout No Access Users"
$ws.Range("C23").Value = "COMPLAINT"
$ws.Range("G23").Value = "mandatory deny read to No Access"
$ws.Rows.Item(23).RowHeight = 23.5

# Row 24 - Case File - Assignee Read Access
$ws.Range("B24").Value = "Case File – Assignee Read Access"
$ws.Range("C24").Value = "CASE_FILE"
$ws.Range("C24").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("G24").Value = "grant read to assignee"
$ws.Rows.Item(24).RowHeight = 23.65

# Row 25 - Case File - Restrict Access to Drafts
$ws.Range("B25").Value = "Case File – Restrict Access to Drafts"
$ws.Range("C25").Value = "CASE_FILE"
$ws.Range("C25").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("D25").Value = "status == 'DRAFT'"
$ws.Range("D25").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("G25").Value = "deny read to *"
$ws.Rows.Item(25).RowHeight = 23.65

# Row 26 - Case File - Grant Access to non-Drafts
$ws.Range("B26").Value = "Case File – Grant Access to non-Drafts"
$ws.Range("C26").Value = "CASE_FILE"
$ws.Range("C26").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("D26").Value = "status != 'DRAFT'"
$ws.Range("D26").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("G26").Value = "grant read to *"
$ws.Rows.Item(26).RowHeight = 23.65

# Row 27 - Case File - Lockout No Access Users
$ws.Range("B27").Value = "Case File – Lockout No Access Users"
$ws.Range("C27").Value = "CASE_FILE"
$ws.Range("G27").Value = "mandatory deny read to No Access"
$ws.Rows.Item(27).RowHeight = 23.5

# New blank row 32 matching the formatting of row 31
$ws.Rows.Item(32).RowHeight = 13.3
$ws.Range("A32").Interior.ColorIndex = $ws.Range("A2").Interior.ColorIndex

Write-Host "done"
